# DEV 14 - Final Fix
# Remove the duplicate/incorrect "Savannah Condopark" project row (row 5).
# This shifts the following row (Archipelago / Bedok, project id 5) up to
# become the new row 5, and shrinks the used range to A1:O5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(5).Delete()
